$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B2').Value = '1.000 (1.000 ± 0.000)'
$ws.Range('C2').Value = '00:04:49 (00:05:20 ± 00:00:23)'
$ws.Range('D2').Value = '00:00:01 (00:00:05 ± 00:00:02)'
$ws.Range('B3').Value = '1.000 (0.997 ± 0.003)'
$ws.Range('C3').Value = '00:00:23 (00:00:25 ± 00:00:01)'
$ws.Range('D3').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B4').Value = '1.000 (1.000 ± 0.001)'
$ws.Range('C4').Value = '00:01:03 (00:01:23 ± 00:00:15)'
$ws.Range('D4').Value = '00:00:01 (00:00:01 ± 00:00:00)'
$ws.Range('B5').Value = '1.000 (1.000 ± 0.000)'
$ws.Range('C5').Value = '00:05:04 (00:05:12 ± 00:00:03)'
$ws.Range('D5').Value = '00:00:01 (00:00:02 ± 00:00:01)'
$ws.Range('B6').Value = '1.000 (1.000 ± 0.000)'
$ws.Range('C6').Value = '00:04:58 (00:05:02 ± 00:00:02)'
$ws.Range('D6').Value = '00:00:01 (00:00:03 ± 00:00:00)'
$ws.Range('B7').Value = '1.000 (1.000 ± 0.000)'
$ws.Range('C7').Value = '00:05:00 (00:05:01 ± 00:00:01)'
$ws.Range('D7').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B8').Value = '1.000 (1.000 ± 0.002)'
$ws.Range('C8').Value = '00:04:45 (00:07:12 ± 00:02:36)'
$ws.Range('D8').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B9').Value = '1.000 (1.000 ± 0.000)'
$ws.Range('C9').Value = '00:04:59 (00:04:59 ± 00:00:00)'
$ws.Range('D9').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B10').Value = '1.000 (0.999 ± 0.002)'
$ws.Range('C10').Value = '00:04:29 (00:04:29 ± 00:00:00)'
$ws.Range('D10').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B11').Value = '0.915 (0.706 ± 0.214)'
$ws.Range('C11').Value = '00:05:05 (00:05:06 ± 00:00:01)'
$ws.Range('D11').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B12').Value = '0.150 (0.104 ± 0.022)'
$ws.Range('C12').Value = '00:02:44 (00:02:51 ± 00:00:04)'
$ws.Range('D12').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B13').Value = '0.601 (0.561 ± 0.027)'
$ws.Range('C13').Value = '00:00:03 (00:00:07 ± 00:00:02)'
$ws.Range('D13').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B14').Value = '1.000 (1.000 ± 0.000)'
$ws.Range('C14').Value = '00:00:56 (00:00:59 ± 00:00:01)'
$ws.Range('D14').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B15').Value = '1.000 (1.000 ± 0.000)'
$ws.Range('C15').Value = '00:01:08 (00:01:14 ± 00:00:08)'
$ws.Range('D15').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B16').Value = '1.000 (0.999 ± 0.002)'
$ws.Range('C16').Value = '00:00:54 (00:00:55 ± 00:00:00)'
$ws.Range('D16').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B17').Value = '1.000 (0.983 ± 0.076)'
$ws.Range('C17').Value = '00:05:01 (00:05:19 ± 00:00:13)'
$ws.Range('D17').Value = '00:00:00 (00:00:00 ± 00:00:00)'
